# Update country data file: "Data" sheet -> "Summary" sheet, with refreshed
# MSME table content (new "Source Type" row, new 86.6 data point, and an
# expanded source citation block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet -----------------------------------------------
$ws.Name = "Summary"

# --- 2. Wipe existing content/formatting so we can lay the sheet out
#        fresh at its new row positions ---------------------------------
$ws.Cells.Clear()

# --- 3. Title ------------------------------------------------------------
$ws.Range("A1").Value = "Guatemala"
$ws.Range("A1").Font.Size = 18

# --- 4. Section header -----------------------------------------------
$ws.Range("A3").Value = "MSME Participation on the Economy"
$ws.Range("A3").Font.Bold = $true

# --- 5. New "Source Type" row (bold + underline) ------------------------
$ws.Range("A9").Value = "Source Type: SME Associations (Most Widely Used)"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Underline = $true

# --- 6. Column headers (row 11) ------------------------------------------
$ws.Range("B11").Value = "Micro"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "SMEs"
$ws.Range("C11").Font.Bold = $true
$ws.Range("D11").Value = "MSMEs"
$ws.Range("D11").Font.Bold = $true

# --- 7. Data rows (12-16): row label bold, values as plain text ---------
$ws.Range("A12").Value = "Enterprises (absolute #)"
$ws.Range("A12").Font.Bold = $true
$ws.Range("B12").Value = "'175566"
$ws.Range("C12").Value = "'8902"
$ws.Range("D12").Value = "'184468"

$ws.Range("A13").Value = "Enterprises density (per 1000 people)"
$ws.Range("A13").Font.Bold = $true
$ws.Range("B13").Value = "'12.9"
$ws.Range("C13").Value = "'0.7"
$ws.Range("D13").Value = "'13.5"

$ws.Range("A14").Value = "Employment (% of total)"
$ws.Range("A14").Font.Bold = $true
$ws.Range("B14").Value = "'68"
$ws.Range("C14").Value = "'18.7"
$ws.Range("D14").Value = "'86.6"

$ws.Range("A15").Value = "Employment (absolute #)"
$ws.Range("A15").Font.Bold = $true
$ws.Range("B15").Value = "'3713900"
$ws.Range("C15").Value = "'1019500"
$ws.Range("D15").Value = "'4733400"

$ws.Range("A16").Value = "Enterprises (% of total)"
$ws.Range("A16").Font.Bold = $true
$ws.Range("B16").Value = "'94.4"
$ws.Range("C16").Value = "'4.8"
$ws.Range("D16").Value = "'99.2"

# --- 8. Source line (italic) ---------------------------------------------
$ws.Range("A17").Value = "Source: CIEN, 2008"
$ws.Range("A17").Font.Italic = $true

# --- 9. New citation block -------------------------------------------
$ws.Range("A23").Value = "CIEN"
$ws.Range("A23").Font.Bold = $true

$ws.Range("A24").Value = "CENTRO DE INVESTIGACIONES ECONÓMICAS NACIONALES, MICRO, PEQUENAS Y MEDIANAS EMPRESAS EN GUATEMALA. Available at http://www.mejoremosguate.org/cms/content/files/diagnosticos/economicos/Lineamientos_PYMES_05-05-2011.pdf"
$ws.Range("A24").Font.Italic = $true
